# The workbook originally contains a small two-column glossary header
# ("単語" / "意味") on Sheet1. This edit appends a new word entry in
# column A, row 2: "こうせいか" (matching the shared-string / sheetData
# content shown in the target OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new word in A2 (unformatted, same as the target cell <c r="A2" t="s">)
$ws.Range("A2").Value = "こうせいか"

# Widen column A to fit the new Japanese text (the saved file shows the
# column auto-sized to roughly 10.27 characters wide with customWidth set).
$ws.Columns.Item(1).ColumnWidth = 9.45

# Match the saved file's last active selection (D9).
[void]$ws.Range("D9").Select()
